$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-9 from 45243 to 45244
# (2023-11-13 -> 2023-11-14), keeping existing cell formatting/style intact.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 45244
}
